$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28; this shifts rows 28..76 down to 29..77
# and updates the sheet dimension automatically.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the latest weekly data point.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0).AddDays(45246)
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112026
$ws.Cells.Item(28, 7).Value = "Haba"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 12).Value = 10000
$ws.Cells.Item(28, 13).Value = 10000
$ws.Cells.Item(28, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 400
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
